$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '53.554.05'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.195.75'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -7.05%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '484.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '125.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -4.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.216.21'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.22%  '
$ws.Range('E10').Value = '  -6.30%  '
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('E12').Value = '  -3.58%  '
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.586.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.96%  '
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '53.479.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000128'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.210.40'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -8.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.34%  '
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '294.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.24%  '
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.145'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.299.88'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '165.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.23%  '
$ws.Range('E31').Value = '  -3.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = '0.0₃0663'
$ws.Range('E33').Value = '  -6.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.994'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.69'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('E37').Value = '  -1.53%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.824'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.92%  '
$ws.Range('E40').Value = '  -4.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '35.74'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.366'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '124.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.76'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('E47').Value = '  -2.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.533'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '232.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0469'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('E51').Value = '  -2.98%  '
